$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Scanner" to "Pharmacology"
$ws.Name = "Pharmacology"

# Convert Log Date column (C2:C13) from MM/DD/YYYY to DD/MM/YYYY text format
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = "20/05/2025"
}
